# The only meaningful content change in this revision is that the
# "s1Protocol" value (column H, shared by every data row) was corrected
# from "E7760" to "E7420". Update the whole data range in one assignment
# so every row keeps pointing at a single shared string, just like the
# original sheet did.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H2:H37").Value = "E7420"
